$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("z")

$data = @(
    @("EU27", "Activity", "Production of electricity by biomass and waste", "EU27", "Commodity", "Electricity", "Update", 0.04375),
    @("EU27", "Activity", "Production of electricity by coal", "EU27", "Commodity", "Electricity", "Update", 0),
    @("EU27", "Activity", "Production of electricity by gas", "EU27", "Commodity", "Electricity", "Update", 0.13125),
    @("EU27", "Activity", "Production of electricity by hydro", "EU27", "Commodity", "Electricity", "Update", 0.08636363636363636),
    @("EU27", "Activity", "Production of electricity by nuclear", "EU27", "Commodity", "Electricity", "Update", 0.1568181818181818),
    @("EU27", "Activity", "Production of electricity by petroleum and other oil derivatives", "EU27", "Commodity", "Electricity", "Update", 0),
    @("EU27", "Activity", "Production of electricity by solar photovoltaic", "EU27", "Commodity", "Electricity", "Update", 0.1454545454545454),
    @("EU27", "Activity", "Production of electricity by wind", "EU27", "Commodity", "Electricity", "Update", 0.4363636363636363)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
